$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "314.74" },
    @{ Cell = "E2"; Value = "3.41%" },
    @{ Cell = "D3"; Value = "35.96" },
    @{ Cell = "E3"; Value = "0.90%" },
    @{ Cell = "D4"; Value = "5.124" },
    @{ Cell = "E4"; Value = "0.81%" },
    @{ Cell = "D5"; Value = "0.08124" },
    @{ Cell = "E5"; Value = "3.29%" },
    @{ Cell = "D6"; Value = "2.122" },
    @{ Cell = "E6"; Value = "-0.20%" },
    @{ Cell = "D7"; Value = "8.022" },
    @{ Cell = "E7"; Value = "1.38%" },
    @{ Cell = "D8"; Value = "0.9294" },
    @{ Cell = "E8"; Value = "1.24%" },
    @{ Cell = "D9"; Value = "0.1018" },
    @{ Cell = "E9"; Value = "4.45%" },
    @{ Cell = "D10"; Value = "0.1876" },
    @{ Cell = "E10"; Value = "1.13%" },
    @{ Cell = "D11"; Value = "0.09159" },
    @{ Cell = "E11"; Value = "6.73%" },
    @{ Cell = "D12"; Value = "0.03587" },
    @{ Cell = "E12"; Value = "1.09%" },
    @{ Cell = "D13"; Value = "0.09922" },
    @{ Cell = "E13"; Value = "-0.20%" },
    @{ Cell = "D14"; Value = "0.001439" },
    @{ Cell = "E14"; Value = "-0.50%" },
    @{ Cell = "D15"; Value = "0.005760" },
    @{ Cell = "E15"; Value = "2.14%" },
    @{ Cell = "D16"; Value = "3.473" },
    @{ Cell = "E16"; Value = "0.29%" },
    @{ Cell = "E17"; Value = "1.15%" },
    @{ Cell = "D18"; Value = "2.704" },
    @{ Cell = "E18"; Value = "5.86%" },
    @{ Cell = "D19"; Value = "0.3370" },
    @{ Cell = "E19"; Value = "-1.62%" },
    @{ Cell = "D20"; Value = "0.1334" },
    @{ Cell = "E20"; Value = "1.90%" },
    @{ Cell = "D21"; Value = "5.135" },
    @{ Cell = "E21"; Value = "-1.72%" },
    @{ Cell = "D23"; Value = "0.04578" },
    @{ Cell = "E23"; Value = "0.71%" },
    @{ Cell = "D24"; Value = "0.001251" },
    @{ Cell = "E24"; Value = "1.29%" },
    @{ Cell = "D25"; Value = "0.004721" },
    @{ Cell = "E25"; Value = "-6.66%" },
    @{ Cell = "D26"; Value = "0.0001255" },
    @{ Cell = "E26"; Value = "-21.73%" },
    @{ Cell = "D27"; Value = "0.0004516" },
    @{ Cell = "E27"; Value = "-4.89%" },
    @{ Cell = "D39"; Value = "0.01983" },
    @{ Cell = "E39"; Value = "7.27%" },
    @{ Cell = "D40"; Value = "0.04912" },
    @{ Cell = "E40"; Value = "4.00%" },
    @{ Cell = "D41"; Value = "0.007857" },
    @{ Cell = "E41"; Value = "4.89%" },
    @{ Cell = "D42"; Value = "0.1394" },
    @{ Cell = "E42"; Value = "-0.14%" },
    @{ Cell = "D43"; Value = "0.007841" },
    @{ Cell = "E43"; Value = "1.24%" },
    @{ Cell = "D44"; Value = "0.002111" },
    @{ Cell = "E44"; Value = "-4.28%" },
    @{ Cell = "D45"; Value = "0.01163" },
    @{ Cell = "E45"; Value = "5.67%" },
    @{ Cell = "D46"; Value = "0.00006513" },
    @{ Cell = "E46"; Value = "2.90%" },
    @{ Cell = "D47"; Value = "0.00000000753" },
    @{ Cell = "E47"; Value = "0.38%" },
    @{ Cell = "D48"; Value = "35.34" },
    @{ Cell = "E48"; Value = "-24.16%" },
    @{ Cell = "D49"; Value = "0.001907" },
    @{ Cell = "E49"; Value = "-4.63%" },
    @{ Cell = "D50"; Value = "0.00002107" },
    @{ Cell = "E50"; Value = "0.38%" },
    @{ Cell = "D51"; Value = "0.0002007" },
    @{ Cell = "E51"; Value = "0.38%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
